$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114 column C is updated to a new translation ("מוצרים נלווים" - "related products"),
# replacing the former "תוספות למוצרים".
$ws.Range("C114").Value = "מוצרים נלווים"

# Two new rows are appended with new English/Hebrew string pairs.
$ws.Range("B115").Value = "All attributes"
$ws.Range("C115").Value = "כל הנלווים"

$ws.Range("B116").Value = "Add attribute"
$ws.Range("C116").Value = "הוסף נלווה"

# Move the active selection to C115 (matches the saved workbook view state).
$ws.Range("C115").Select()
